# Generate Report for Handoff
#
# Refreshing the "Ready for handoff" rows: the handoff-xliff-generation
# timestamps move forward a few seconds, and the Priority column for
# those rows now records the handoff type ("ht").

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

# Rows whose source file is in "Ready for handoff" status (row 12 -
# 9d6903c6-... - is excluded, matching the source data).
$rows = @(7, 8, 9, 10, 11, 13)

foreach ($r in $rows) {
    # Overview!G<r>: Latest HO Xliff Generate Date
    $overview.Range("G$r").Value = "2016-08-16 10:21:26"

    # zh-cn!H<r>: Latest Handoff Datetime
    $zhcn.Range("H$r").Value = "2016-08-16 10:21:20"

    # de-de!H<r>: Latest Handoff Datetime
    $dede.Range("H$r").Value = "2016-08-16 10:21:26"

    # Priority column now set to "ht" for both locale sheets
    $zhcn.Range("E$r").Value = "ht"
    $dede.Range("E$r").Value = "ht"
}
